$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right count 4 -> 5, Wrong-answer penalty -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 ("Total"): Right total 28 -> 35, deduction total 0 -> -0 (0 wrong * -1.2),
# and the score summary text 28/112 -> 35.0/140
$ws.Range("B12").Value = 35
$ws.Range("C12").Value = -0.0
$ws.Range("E12").Value = "35.0/140"
